$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 4), mirroring row 3's content/format but with
# the next sequential NUM value (validaciones en carga masiva -> nuevo
# registro de prueba).
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = $ws.Range("B3").Value2
$ws.Range("C4").Value = $ws.Range("C3").Value2
$ws.Range("D4").Value = 665544
$ws.Range("E4").Value = "B@SRI.COM"
$ws.Range("F4").Value = 67890
$ws.Range("G4").Value = $ws.Range("G3").Value2
$ws.Range("H4").Value = $ws.Range("H3").Value2

# E4 carries the same hyperlink/style treatment as E2/E3 (mailto link,
# "Hipervinculo" cell style).
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:B@SRI.COM")
$ws.Range("E4").Style = $ws.Range("E3").Style

$ws.Range("A4").Select()
